# Update "paises.xlsx" (sheet "Pais") with a refreshed COVID-19 data pull
# (20:22 -> 20:52) and the resulting re-sort of several rows by "Casos totales".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 20:52"

# --- Rows whose country swapped places with its neighbour after the re-sort ---
# (India / Israel)
$ws.Range("A22").Value = "Israel"
$ws.Range("B22").Value = 12501
$ws.Range("C22").Value = 455
$ws.Range("D22").Value = 2563
$ws.Range("E22").Value = 9808
$ws.Range("F22").Value = 180
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 130

$ws.Range("A23").Value = "India"
$ws.Range("B23").Value = 12322
$ws.Range("C23").Value = 835
$ws.Range("D23").Value = 1432
$ws.Range("E23").Value = 10485
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 405

# (Corea del Sur / Peru)
$ws.Range("A25").Value = "Peru"
$ws.Range("B25").Value = 11475
$ws.Range("C25").Value = 1172
$ws.Range("D25").Value = 3108
$ws.Range("E25").Value = 8113
$ws.Range("F25").Value = 146
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 254

$ws.Range("A26").Value = "Corea del Sur"
$ws.Range("B26").Value = 10591
$ws.Range("C26").Value = 27
$ws.Range("D26").Value = 7616
$ws.Range("E26").Value = 2750
$ws.Range("F26").Value = 55
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 225

# --- Rows with an updated case count only (no reordering) ---
$ws.Range("E31").Value = 5627
$ws.Range("G31").Value = 21
$ws.Range("H31").Value = 372

$ws.Range("D90").Value = 85
$ws.Range("E90").Value = 552

$ws.Range("B117").Value = 238
$ws.Range("C117").Value = 5
$ws.Range("E117").Value = 168

$ws.Range("B130").Value = 131
$ws.Range("C130").Value = 2
$ws.Range("E130").Value = 27

# --- "Somalia" moves up next to "Gabon", shifting Liechtenstein/Birmania/Barbados down one row ---
$ws.Range("A143").Value = "Somalia"
$ws.Range("B143").Value = 80
$ws.Range("C143").Value = 20
$ws.Range("D143").Value = 2
$ws.Range("E143").Value = 73
$ws.Range("F143").Value = 2
$ws.Range("G143").Value = 3
$ws.Range("H143").Value = 5

$ws.Range("A144").Value = "Liechtenstein"
$ws.Range("B144").Value = 79
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 55
$ws.Range("E144").Value = 23
$ws.Range("H144").Value = 1

$ws.Range("A145").Value = "Birmania"
$ws.Range("B145").Value = 74
$ws.Range("C145").Value = 11
$ws.Range("D145").Value = 2
$ws.Range("E145").Value = 68
$ws.Range("F145").Value = 0
$ws.Range("H145").Value = 4

$ws.Range("A146").Value = "Barbados"
$ws.Range("B146").Value = 73
$ws.Range("D146").Value = 15
$ws.Range("E146").Value = 53
$ws.Range("F146").Value = 4
$ws.Range("H146").Value = 5
